# Remove the "Prayer Requests (Image)" slide (position 2). The hymnal /
# prayer-request detection no longer relies on a dedicated image slide with
# this title; the remaining slides shift up to fill the gap automatically.
$p = $ppt.ActivePresentation

$targetIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "Prayer Requests (Image)") {
                $targetIndex = $i
                break
            }
        }
    }
    if ($targetIndex -ne -1) { break }
}

if ($targetIndex -eq -1) {
    # Fallback: the slide is expected to be the second slide in the deck.
    $targetIndex = 2
}

$p.Slides.Item($targetIndex).Delete()
